$d = $word.ActiveDocument

$d.Content.Find.Execute("## 1   23.41360 -2.466498", $true, $false, $false, $false, $false, $true, 1, $false, "## 1   23.00168 -2.466498", 2)
$d.Content.Find.Execute("## 2   23.31524 -2.564851", $true, $false, $false, $false, $false, $true, 1, $false, "## 2   22.90333 -2.564851", 2)
$d.Content.Find.Execute("## 3   24.59679 -1.283304", $true, $false, $false, $false, $false, $true, 1, $false, "## 3   24.18488 -1.283304", 2)
$d.Content.Find.Execute("## 4   26.53007  0.649978", $true, $false, $false, $false, $false, $true, 1, $false, "## 4   26.11816  0.649978", 2)
$d.Content.Find.Execute("## 5   24.16439 -1.715703", $true, $false, $false, $false, $false, $true, 1, $false, "## 5   23.75248 -1.715703", 2)
$d.Content.Find.Execute("## 6   24.14983 -1.730263", $true, $false, $false, $false, $false, $true, 1, $false, "## 6   23.73792 -1.730263", 2)

$d.Content.Find.Execute("##     21.09169     22.87938     24.51027     26.04739", $true, $false, $false, $false, $false, $true, 1, $false, "##     20.67977     22.46747     24.09836     25.63548", 2)
